$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 2.63
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 2.03
$ws.Range("S2").Value = 2.4
$ws.Range("T2").Value = 1.53
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 1.18
$ws.Range("Y2").Value = 1.57
$ws.Range("AC2").Value = 7
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 29
$ws.Range("AN2").Value = 7
$ws.Range("AR2").Value = 26
